# Final tweaks to l16: add a "click to build by paragraph" entrance
# animation (Appear) to the two content placeholders on the
# "incremental content" slide, so each of their 5 paragraphs appears
# one click at a time.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$seq = $s.TimeLine.MainSequence

# Shape.Id=3 "Content Placeholder 2" and Shape.Id=4 "Content Placeholder 3"
# are Shapes.Item(2) and Shapes.Item(3) respectively on this slide.
$leftBox = $s.Shapes.Item(2)
$rightBox = $s.Shapes.Item(3)

# effectId 1 = ppEffectAppear (entrance "Appear"); Level 2 = build by
# 1st-level paragraph, which expands into one click-effect per paragraph.
$leftEffect = $seq.AddEffect($leftBox, 1, 2)
$rightEffect = $seq.AddEffect($rightBox, 1, 2)
